$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Modify existing rows 92 and 93 ---
$ws.Range("F92").Value = 20163126.5
$ws.Range("G92").Value = 500
$ws.Range("H92").Value = 469

$ws.Range("F93").Value = 57805127.3
$ws.Range("G93").Value = 300
$ws.Range("H93").Value = 972

# --- Append new rows 577-587 ---
# row 577
$ws.Range("A577").NumberFormat = "@"
$ws.Range("A577").Value = "2024-05-14"
$ws.Range("A577").Style = "Normal"
$ws.Range("B577").Value = 869
$ws.Range("C577").Value = 884
$ws.Range("D577").Value = 895
$ws.Range("E577").Value = 856
$ws.Range("F577").Value = 951634104.5
$ws.Range("G577").Value = 1084400
$ws.Range("H577").Value = 3274

# row 578
$ws.Range("A578").NumberFormat = "@"
$ws.Range("A578").Value = "2024-05-15"
$ws.Range("A578").Style = "Normal"
$ws.Range("B578").Value = 902.5
$ws.Range("C578").Value = 875
$ws.Range("D578").Value = 907
$ws.Range("E578").Value = 870
$ws.Range("F578").Value = 843252407
$ws.Range("G578").Value = 939978
$ws.Range("H578").Value = 3445

# row 579
$ws.Range("A579").NumberFormat = "@"
$ws.Range("A579").Value = "2024-05-16"
$ws.Range("A579").Style = "Normal"
$ws.Range("B579").Value = 919
$ws.Range("C579").Value = 903
$ws.Range("D579").Value = 923.5
$ws.Range("E579").Value = 890
$ws.Range("F579").Value = 984813526
$ws.Range("G579").Value = 1080910
$ws.Range("H579").Value = 3515

# row 580
$ws.Range("A580").NumberFormat = "@"
$ws.Range("A580").Value = "2024-05-17"
$ws.Range("A580").Style = "Normal"
$ws.Range("B580").Value = 938
$ws.Range("C580").Value = 918
$ws.Range("D580").Value = 949.5
$ws.Range("E580").Value = 914.5
$ws.Range("F580").Value = 1182251724.5
$ws.Range("G580").Value = 1265485
$ws.Range("H580").Value = 3965

# row 581
$ws.Range("A581").NumberFormat = "@"
$ws.Range("A581").Value = "2024-05-20"
$ws.Range("A581").Style = "Normal"
$ws.Range("B581").Value = 966
$ws.Range("C581").Value = 937
$ws.Range("D581").Value = 972
$ws.Range("E581").Value = 925
$ws.Range("F581").Value = 1019462117.5
$ws.Range("G581").Value = 1078830
$ws.Range("H581").Value = 6604

# row 582
$ws.Range("A582").NumberFormat = "@"
$ws.Range("A582").Value = "2024-05-21"
$ws.Range("A582").Style = "Normal"
$ws.Range("B582").Value = 1030
$ws.Range("C582").Value = 972
$ws.Range("D582").Value = 1031
$ws.Range("E582").Value = 972
$ws.Range("F582").Value = 1889732348.5
$ws.Range("G582").Value = 1878198
$ws.Range("H582").Value = 5948

# row 583
$ws.Range("A583").NumberFormat = "@"
$ws.Range("A583").Value = "2024-05-22"
$ws.Range("A583").Style = "Normal"
$ws.Range("B583").Value = 1050
$ws.Range("C583").Value = 1035
$ws.Range("D583").Value = 1053
$ws.Range("E583").Value = 1018
$ws.Range("F583").Value = 1038857078
$ws.Range("G583").Value = 998395
$ws.Range("H583").Value = 8863

# row 584
$ws.Range("A584").NumberFormat = "@"
$ws.Range("A584").Value = "2024-05-23"
$ws.Range("A584").Style = "Normal"
$ws.Range("B584").Value = 1015.5
$ws.Range("C584").Value = 1055
$ws.Range("D584").Value = 1055
$ws.Range("E584").Value = 1000
$ws.Range("F584").Value = 886611054.5
$ws.Range("G584").Value = 871269
$ws.Range("H584").Value = 3922

# row 585
$ws.Range("A585").NumberFormat = "@"
$ws.Range("A585").Value = "2024-05-24"
$ws.Range("A585").Style = "Normal"
$ws.Range("B585").Value = 1015
$ws.Range("C585").Value = 1015.5
$ws.Range("D585").Value = 1028.5
$ws.Range("E585").Value = 980
$ws.Range("F585").Value = 1303010146.5
$ws.Range("G585").Value = 1303984
$ws.Range("H585").Value = 3624

# row 586
$ws.Range("A586").NumberFormat = "@"
$ws.Range("A586").Value = "2024-05-27"
$ws.Range("A586").Style = "Normal"
$ws.Range("B586").Value = 1002.5
$ws.Range("C586").Value = 990
$ws.Range("D586").Value = 1019.5
$ws.Range("E586").Value = 990
$ws.Range("F586").Value = 194485919
$ws.Range("G586").Value = 1
$ws.Range("H586").Value = 1228

# row 587
$ws.Range("A587").NumberFormat = "@"
$ws.Range("A587").Value = "2024-05-28"
$ws.Range("A587").Style = "Normal"
$ws.Range("B587").Value = 1029.5
$ws.Range("C587").Value = 1022
$ws.Range("D587").Value = 1035
$ws.Range("E587").Value = 991.5
$ws.Range("F587").Value = 955303045.5
$ws.Range("G587").Value = 939249
$ws.Range("H587").Value = 2462
